$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, [string]$text) {
    $escaped = $text -replace '"', '""'
    $ws.Range($cellRef).Formula = '="' + $escaped + '"'
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

Set-TextCell "D2" "63.033.07"
Set-TextCell "E2" "  -0.97%  "

Set-TextCell "D3" "2.548.95"
Set-TextCell "E3" "  -0.14%  "

Set-TextCell "E4" "  +0.05%  "

Set-TextCell "D5" "583.05"
Set-TextCell "E5" "  +2.01%  "

Set-TextCell "D6" "146.84"
Set-TextCell "E6" "  -3.00%  "

Set-TextCell "E8" "  -0.58%  "

Set-TextCell "E9" "  -0.61%  "

Set-TextCell "E10" "  -3.61%  "

Set-TextCell "E11" "  -0.12%  "

Set-TextCell "E12" "  -1.52%  "

Set-TextCell "E13" "  -3.55%  "

Set-TextCell "D14" "3.006.18"
Set-TextCell "E14" "  -0.01%  "

Set-TextCell "D15" "62.927.27"
Set-TextCell "E15" "  -0.96%  "

Set-TextCell "E16" "  -1.26%  "

Set-TextCell "E17" "  -0.04%  "

Set-TextCell "E18" "  -3.22%  "

Set-TextCell "D19" "338.39"
Set-TextCell "E19" "  -0.87%  "

Set-TextCell "E20" "  -1.34%  "

Set-TextCell "D21" "6.76"
Set-TextCell "E21" "  -1.59%  "

Set-TextCell "E22" "  -0.34%  "

Set-TextCell "D23" "65.66"
Set-TextCell "E23" "  -0.93%  "

Set-TextCell "D24" "2.678.64"
Set-TextCell "E24" "  +0.37%  "

Set-TextCell "E25" "  -0.85%  "

Set-TextCell "E26" "  -0.51%  "

Set-TextCell "B27" "SuiNetwork"
Set-TextCell "C27" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextCell "D27" "1.49"
Set-TextCell "E27" "  -1.68%  "

Set-TextCell "B28" "Binance-PegBSC-USD"
Set-TextCell "C28" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextCell "D28" "1.00"
Set-TextCell "E28" "  +0.09%  "

Set-TextCell "E29" "  -3.68%  "

Set-TextCell "D30" "7.69"
Set-TextCell "E30" "  +6.42%  "

Set-TextCell "E31" "  +4.41%  "

Set-TextCell "D32" "0.0₃0816"
Set-TextCell "E32" "  -2.20%  "

Set-TextCell "D33" "177.99"
Set-TextCell "E33" "  -0.19%  "

Set-TextCell "E34" "  -2.64%  "

Set-TextCell "D35" "417.43"
Set-TextCell "E35" "  -1.31%  "

Set-TextCell "E36" "  -1.66%  "

Set-TextCell "D37" "19.11"
Set-TextCell "E37" "  -0.58%  "

Set-TextCell "E38" "  +0.01%  "

Set-TextCell "E39" "  -2.71%  "

Set-TextCell "D40" "1.74"
Set-TextCell "E40" "  -2.62%  "

Set-TextCell "E41" "  +0.03%  "

Set-TextCell "D42" "39.78"
Set-TextCell "E42" "  +0.01%  "

Set-TextCell "D43" "150.96"

Set-TextCell "E44" "  -1.32%  "

Set-TextCell "E45" "  -2.06%  "

Set-TextCell "D46" "0.0539"
Set-TextCell "E46" "  +1.32%  "

Set-TextCell "E47" "  -1.53%  "

Set-TextCell "E48" "  -0.12%  "

Set-TextCell "E49" "  -1.00%  "

Set-TextCell "E50" "  -2.28%  "

Set-TextCell "E51" "  -6.71%  "
